# Auto-generated cell updates based on the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.259.97'
$ws.Range('E2').Value = '  -0.91%  '
$ws.Range('D3').Value = '2.647.36'
$ws.Range('E3').Value = '  +2.41%  '
$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E4').Value = '  +0.23%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '606.47'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  +3.19%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '144.40'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.68%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').Value = '2.647.13'
$ws.Range('E9').Value = '  +2.49%  '
$ws.Range('E10').Value = '  +1.12%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '5.63'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('E12').Value = '  +0.29%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.364'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  +3.57%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '27.34'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '3.120.69'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').Value = '63.127.30'
$ws.Range('E16').Value = '  -0.84%  '
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '2.648.42'
$ws.Range('E18').Value = '  +2.73%  '
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '4.44'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  +2.66%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '342.70'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  -0.21%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '6.84'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +2.77%  '
$ws.Range('E23').Value = '  -0.13%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '67.15'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -1.27%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '1.64'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('E26').Value = '  -2.94%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '8.66'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  +4.60%  '
$ws.Range('E28').Value = '  +0.09%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '547.17'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  +14.12%  '
$ws.Range('E30').Value = '  +0.26%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '7.91'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -0.98%  '
$ws.Range('E32').Value = '  +4.29%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '1.79'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +5.66%  '
$ws.Range('D34').Value = '0.0₃0810'
$ws.Range('E34').Value = '  +0.85%  '
$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '172.30'
$cell.Style = 'Normal'
$ws.Range('E35').Value = '  -2.31%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '5.09'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  +11.26%  '
$ws.Range('E37').Value = '  +1.92%  '
$ws.Range('E38').Value = '  +0.08%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '19.12'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +1.09%  '
$ws.Range('E40').Value = '  +5.72%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '171.79'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +7.04%  '
$ws.Range('E42').Value = '  +0.11%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '3.75'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  +0.91%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '22.34'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  +2.71%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.0575'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  +6.60%  '
$ws.Range('E46').Value = '  -0.34%  '
$ws.Range('E47').Value = '  +0.96%  '
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.0961'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -0.34%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '18.84'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +4.00%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.74'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +1.92%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '11.22'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -1.33%  '
